$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 109 (existing data shifts down to
# rows 110-120, which already matches the target state for those rows).
$ws.Rows(109).Insert()

# Populate the newly inserted row 109 with the new weekly record.
$ws.Range("A109").Value = 11
$ws.Range("B109").Value = "Vega Monumental Concepción"
$ws.Range("C109").Value = "Bíobío"
$ws.Range("D109").Value = 44769
$ws.Range("E109").Value = 8
$ws.Range("F109").Value = 100112021
$ws.Range("G109").Value = "Ají"
$ws.Range("H109").Value = "Inferno"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 100
$ws.Range("K109").Value = 16000
$ws.Range("L109").Value = 17000
$ws.Range("M109").Value = 16500
$ws.Range("N109").Value = "$/caja 12 kilos"
$ws.Range("O109").Value = "Región de Arica y Parinacota"
$ws.Range("P109").Value = 1375
$ws.Range("Q109").Value = 12
$ws.Range("R109").Value = "Hortaliza"
